# "include Tanix TX3 mini"
# Insert a new benchmark entry for the "Tanix TX3 mini" (Amlogic S905W)
# device into the "2024" results table, keeping it sorted by the
# Integer Index (column E) the same way the rest of the sheet is sorted.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# The table is sorted ascending by column E ("Integer Index"). The new
# device's Integer Index (10.948) belongs right after row 18 (Nexus 4,
# 9.895) and before the old row 19 (Silentium II, 11.927), so insert a
# fresh row at position 19 and shift everything below it down.
$ws.Rows.Item(19).Insert()

$ws.Cells.Item(19, 1).Value = "Tanix TX3 mini"
$ws.Cells.Item(19, 2).Value = "Amlogic S905W"
$ws.Cells.Item(19, 3).Value = 1200
$ws.Cells.Item(19, 4).Value = 8.68
$ws.Cells.Item(19, 5).Value = 10.948
$ws.Cells.Item(19, 6).Value = 10.155
$ws.Cells.Item(19, 7).Value = 39.718
$ws.Cells.Item(19, 8).Value = 18.31
$ws.Cells.Item(19, 9).Value = "-"
$ws.Cells.Item(19, 10).Value = 45319

# Restore the view to what's left after the edit: scrolled back near the
# top of the frozen pane, with F17 selected.
$ws.Range("F17").Select()
$excel.ActiveWindow.ScrollRow = 2
